# Update Release-Notes.xlsx - Folder inventory updated on Thu Jun 12 06:48:45 UTC 2025

$wb = $excel.ActiveWorkbook

# --- Sheet: Folder Inventory ---
$wsInv = $wb.Worksheets.Item("Folder Inventory")

# The "Cloud-Native Applications" folder was refreshed (new timestamp) and
# is now the most recently updated folder, so it moves from row 11 to the
# top of the data (row 2). Remove its stale row further down the list and
# insert a new row right after the header with the refreshed data.
$wsInv.Rows.Item(11).Delete()
$wsInv.Rows.Item(2).Insert()
$wsInv.Range("A2:E2").ClearFormats()

$wsInv.Range("A2").Value = "Cloud-Native Applications"
$wsInv.Range("B2").Value = "Cloud-Native Applications"
$wsInv.Range("C2").Value = "2025-06-12 12:18:28 +0530"
$wsInv.Range("D2").Value = 1
$wsInv.Range("E2").Value = "Root"

# --- Sheet: Metadata ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2025-06-12 06:48:45 UTC"

# "Workflow Run" is stored as text, not a number; force text formatting so the
# numeric-looking value "7" isn't auto-converted to a numeric cell, then drop
# the temporary number format so no stray style is left behind on the cell.
$wsMeta.Range("B5").NumberFormat = "@"
$wsMeta.Range("B5").Value = "7"
$wsMeta.Range("B5").ClearFormats()

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 12:18:28 +0530"
